# Fixed bug when calculating single worksheets or ranges
#
# Adds a new "If" example block to the ValidateFormulas worksheet (the
# workbook's internal file name is sheet3.xml). A new row is inserted
# above the existing "Sum Kolumn" / "Boolean" example block (which used
# to be row 31 and becomes row 32), and the new row 31 is populated with
# three IF() formula examples plus a row label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidateFormulas")

# Insert a new blank row at position 31; everything below (old rows
# 31-41) shifts down by one (new rows 32-42), and all formula
# references / merged cells are automatically adjusted by Excel.
$ws.Rows.Item(31).Insert()

# Row label in column A, using the same bold header style as the other
# row labels in this block (e.g. "Sum Kolumn" in B32).
$ws.Range("A31").Value = "If"
$ws.Range("A31").Font.Bold = $true

# New IF() formula examples.
$ws.Range("B31").Formula = "=IF(B2>3,B3,B5)"
$ws.Range("C31").Formula = "=IF((B2*B3)*C1<0,(B2*B3)*C1,ABS((B2*B3)*C1))"
$ws.Range("D31").Formula = "=IF((B2*B3)*C1<0,ABS((B2*B3)*C1),(B2*B3)*C1)"

# Match the workbook's recorded selection after the edit.
$ws.Range("D31").Select()
